$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data region (rows 2-16) since we are shrinking the table
$ws.Range("A2:D16").ClearContents()

# New term labels (habitat_type, season, day_night and their interactions - "length" terms removed)
$terms = @(
    "habitat_type",
    "season",
    "day_night",
    "habitat_type:season",
    "habitat_type:day_night",
    "season:day_night",
    "habitat_type:season:day_night"
)

$statistic = @(
    207.133081056651,
    44.9557804044391,
    1667.96395357593,
    268.567726862871,
    149.141553716142,
    144.247607235722,
    122.286658891661
)

$df = @(
    4,
    3,
    3,
    12,
    12,
    9,
    35
)

$pvalue = @(
    0.000000000000000000000000000000000000000000109904551627884,
    0.000000000945510096305254,
    0,
    0.0000000000000000000000000000000000000000000000000181362979053718,
    0.0000000000000000000000000846626344293372,
    0.0000000000000000000000000136757646310861,
    0.0000000000132162470534377
)

for ($i = 0; $i -lt $terms.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $terms[$i]
    $ws.Cells.Item($row, 2).Value = $statistic[$i]
    $ws.Cells.Item($row, 3).Value = $df[$i]
    $ws.Cells.Item($row, 4).Value = $pvalue[$i]
}
